# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计", holding the new
#    quarter's fund-holding detail rows (same shape/style as the other
#    quarterly sheets -- built by duplicating "2022-Q2" so borders/fonts on
#    the index column and header row match exactly, then overwriting every
#    cell with the 2022-Q3 figures).
# 2. Insert a new row at the top of the "总计" (summary) sheet's data table
#    with the 2022-Q3 totals, shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert new "2022-Q3" sheet right after "总计", cloned from "2022-Q2"
#    so it starts out with the same column widths / header & index-column
#    styling, then gets its data replaced wholesale.
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $zongji)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Data rows (fund code, fund name, fund size, total equity position,
# position ratio, held market value, position rank)
$data = @(
  @("000566", "华泰柏瑞创新升级混合A",             "18.10", "89.06", "4.42", "0.8000", 2),
  @("014292", "嘉实产业领先混合A",                 "13.38", "91.87", "5.39", "0.7212", 1),
  @("007968", "华泰柏瑞研究精选混合A",             "8.52",  "88.08", "4.43", "0.3774", 3),
  @("009636", "华泰柏瑞景气优选混合A",             "8.26",  "90.03", "4.42", "0.3651", 3),
  @("000967", "华泰柏瑞创新动力灵活配置混合",       "5.45",  "89.89", "4.42", "0.2409", 3),
  @("013431", "华泰柏瑞景气汇选三年持有期混合A",     "4.89",  "88.55", "4.41", "0.2156", 3),
  @("008373", "华泰柏瑞景气回报一年持有期混合A",     "2.49",  "89.50", "4.42", "0.1101", 3),
  @("013847", "华泰柏瑞匠心汇选混合A",             "1.95",  "89.24", "4.41", "0.0860", 3),
  @("013432", "华泰柏瑞景气汇选三年持有期混合C",     "1.29",  "88.55", "4.41", "0.0569", 3),
  @("010291", "华泰柏瑞研究精选混合C",             "1.03",  "88.08", "4.43", "0.0456", 3),
  @("010028", "华泰柏瑞创新升级混合C",             "0.49",  "89.06", "4.42", "0.0217", 2),
  @("011454", "华泰柏瑞景气优选混合C",             "0.43",  "90.03", "4.42", "0.0190", 3),
  @("014293", "嘉实产业领先混合C",                 "0.26",  "91.87", "5.39", "0.0140", 1),
  @("013848", "华泰柏瑞匠心汇选混合C",             "0.18",  "89.24", "4.41", "0.0079", 3),
  @("008374", "华泰柏瑞景气回报一年持有期混合C",     "0.07",  "89.50", "4.42", "0.0031", 3)
)

$rowCount = $data.Count
$lastRow = $rowCount + 1

# Column A: numeric row index (0-based), inherits the template's bold /
# bordered / centered style already present on the copied sheet.
for ($i = 0; $i -lt $rowCount; $i++) {
  $q3.Cells.Item($i + 2, 1).Value2 = [double]$i
}

# Columns B..G: text-like values (fund code/name/size/position/ratio/value)
# -- force text format so numeric-looking strings ("000566", "18.10") are
# not reinterpreted as numbers, matching the template's own data cells.
$bgRange = $q3.Range("B2:G$lastRow")
$bgRange.NumberFormat = "@"
for ($i = 0; $i -lt $rowCount; $i++) {
  $r = $i + 2
  for ($j = 0; $j -lt 6; $j++) {
    $q3.Cells.Item($r, 2 + $j).Value2 = [string]$data[$i][$j]
  }
}

# Column H: numeric rank
for ($i = 0; $i -lt $rowCount; $i++) {
  $q3.Cells.Item($i + 2, 8).Value2 = [double]$data[$i][6]
}

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q3 summary row at the top of the "总计" sheet table
# ---------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

# Match the inserted row's formatting to its neighbours (index column
# keeps the bold/bordered/centered style already used by every other row).
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)

$zongji.Cells.Item(2, 1).Value2 = 0
$zongji.Cells.Item(2, 2).Value2 = "2022-Q3"
$zongji.Cells.Item(2, 3).Value2 = 15
$zongji.Cells.Item(2, 4).Value2 = 3.08

# Re-number the index column (A), since Rows.Insert shifted the existing
# quarters down by one row but left their literal index values untouched.
$zongji.Cells.Item(3, 1).Value2 = 1
$zongji.Cells.Item(4, 1).Value2 = 2
$zongji.Cells.Item(5, 1).Value2 = 3
$zongji.Cells.Item(6, 1).Value2 = 4
$zongji.Cells.Item(7, 1).Value2 = 5
$zongji.Cells.Item(8, 1).Value2 = 6

$zongji.Range("A1").Select()
$excel.CutCopyMode = $false
